$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: now describes PM_IPA_FERMENTACION_PRESION (like rows 2-4), with new
# measurement values and a "presion alta/normal" formula (previously it held
# the PM_IPA_CENTRIFUGADO_MARCHA "green/blue" comparison that now lives on row 6).
$ws.Range("A5").Value = "PM_IPA_FERMENTACION_PRESION"
$ws.Range("B5").Value = 7.45
$ws.Range("C5").Value = 6.8
$ws.Range("D5").Formula = '=IF(AND(B5>3,B5<7),"presion alta","presion normal")'
$ws.Range("F5").Value = "text"

# Row 6: now carries the "green/blue" equality check that row 5 used to have.
$ws.Range("D6").Formula = '=IF(B6=C6,"green","blue")'

# Row 7: now checks for a combined "marcha" state instead of pressure stability.
$ws.Range("D7").Formula = '=IF(AND(B7>3,B7<7),"marcha combinada","marcha no combinada")'

# Row 8 is removed entirely (not just emptied) - clear formatting too so the
# row disappears from the saved sheet instead of lingering as a blank row.
# Rows below row 8 (e.g. row 11) stay put - this is a full clear, not a
# row delete/shift.
$ws.Rows(8).Clear()

# Move the active selection to F2 (was D9).
$ws.Range("F2").Select()
